$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.fiercebiotech.com/medtech/roche-receives-fda-breakthrough-label-ai-powered-lung-cancer-companion-diagnostic-test"
$keyword = "companion diagnostic"
$titleHtml = '<a href="https://www.fiercebiotech.com/medtech/roche-receives-fda-breakthrough-label-ai-powered-lung-cancer-companion-diagnostic-test" hreflang="en">Roche receives FDA breakthrough label for AI-powered lung cancer companion diagnostic test</a>'

# Add a new row (row 6) with link/keywords/title data, mirroring the
# existing rows produced by the feed-filtering workflow.
$ws.Range("A6").Value = $url
$ws.Hyperlinks.Add($ws.Range("A6"), $url) | Out-Null

# Match the "Hyperlink" cell style already used by the other link cells
# (A2:A5) by copying their formatting onto the new cell.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B6").Value = $keyword
$ws.Range("C6").Value = $titleHtml
